$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 2526757.8
$ws.Cells.Item(18, 9).Value = 3087281.8
$ws.Cells.Item(18, 10).Value = 4400
$ws.Cells.Item(18, 11).Value = 3087281.8
$ws.Cells.Item(18, 12).Value = 4400
$ws.Cells.Item(18, 13).Value = -3086997.8
$ws.Cells.Item(18, 14).Value = -4968
$ws.Cells.Item(74, 8).Value = 3734.4285
$ws.Cells.Item(74, 9).Value = 3527.4285
$ws.Cells.Item(74, 10).Value = 3941.4285
$ws.Cells.Item(74, 11).Value = 3527.4285
$ws.Cells.Item(74, 12).Value = 3941.4285
$ws.Cells.Item(74, 13).Value = -2591.4285
$ws.Cells.Item(74, 14).Value = -5813.4285
$ws.Cells.Item(77, 8).Value = 3734.4285
$ws.Cells.Item(77, 9).Value = 3527.4285
$ws.Cells.Item(77, 10).Value = 3941.4285
$ws.Cells.Item(77, 11).Value = 17637.1425
$ws.Cells.Item(77, 12).Value = 19707.1425
$ws.Cells.Item(77, 13).Value = -12957.1425
$ws.Cells.Item(77, 14).Value = -29067.1425
$ws.Cells.Item(113, 8).Value = 3266.389
$ws.Cells.Item(113, 9).Value = 3119.9
$ws.Cells.Item(113, 10).Value = 3449.5
$ws.Cells.Item(113, 11).Value = 3119.9
$ws.Cells.Item(113, 12).Value = 3449.5
$ws.Cells.Item(113, 13).Value = 134.0999999999999
$ws.Cells.Item(113, 14).Value = -9957.5
$ws.Cells.Item(132, 8).Value = 5003020
$ws.Cells.Item(132, 9).Value = 6063015.5
$ws.Cells.Item(132, 11).Value = 18189046.5
$ws.Cells.Item(132, 13).Value = -18186516.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 1596.6666
$ws.Cells.Item(4, 9).Value = 1170
$ws.Cells.Item(4, 10).Value = 2450
$ws.Cells.Item(4, 11).Value = 1170
$ws.Cells.Item(4, 12).Value = 2450
$ws.Cells.Item(4, 13).Value = -1054
$ws.Cells.Item(4, 14).Value = -2682
$ws.Cells.Item(5, 8).Value = 130.4
$ws.Cells.Item(5, 9).Value = 150
$ws.Cells.Item(5, 10).Value = 117.333336
$ws.Cells.Item(5, 11).Value = 150
$ws.Cells.Item(5, 12).Value = 117.333336
$ws.Cells.Item(5, 13).Value = -38
$ws.Cells.Item(5, 14).Value = -341.333336
$ws.Cells.Item(17, 8).Value = 10009
$ws.Cells.Item(17, 10).Value = 10009
$ws.Cells.Item(17, 12).Value = 10009
$ws.Cells.Item(17, 14).Value = -10355
$ws.Cells.Item(32, 8).Value = 5541.85
$ws.Cells.Item(32, 9).Value = 5490.9443
$ws.Cells.Item(32, 10).Value = 6000
$ws.Cells.Item(32, 11).Value = 5490.9443
$ws.Cells.Item(32, 12).Value = 6000
$ws.Cells.Item(32, 13).Value = -5203.9443
$ws.Cells.Item(32, 14).Value = -6574
$ws.Cells.Item(88, 8).Value = 4751.75
$ws.Cells.Item(88, 9).Value = 3000
$ws.Cells.Item(88, 10).Value = 5335.6665
$ws.Cells.Item(88, 11).Value = 3000
$ws.Cells.Item(88, 12).Value = 5335.6665
$ws.Cells.Item(88, 13).Value = -2594
$ws.Cells.Item(88, 14).Value = -6147.6665
$ws.Cells.Item(91, 8).Value = 4751.75
$ws.Cells.Item(91, 9).Value = 3000
$ws.Cells.Item(91, 10).Value = 5335.6665
$ws.Cells.Item(91, 11).Value = 3000
$ws.Cells.Item(91, 12).Value = 5335.6665
$ws.Cells.Item(91, 13).Value = -1596
$ws.Cells.Item(91, 14).Value = -8143.6665
$ws.Cells.Item(97, 8).Value = 638.7143
$ws.Cells.Item(97, 9).Value = 576.75
$ws.Cells.Item(97, 10).Value = 1010.5
$ws.Cells.Item(97, 11).Value = 576.75
$ws.Cells.Item(97, 12).Value = 1010.5
$ws.Cells.Item(97, 13).Value = -80.75
$ws.Cells.Item(97, 14).Value = -2002.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 130.4
$ws.Cells.Item(4, 9).Value = 150
$ws.Cells.Item(4, 10).Value = 117.333336
$ws.Cells.Item(4, 11).Value = 150
$ws.Cells.Item(4, 12).Value = 117.333336
$ws.Cells.Item(4, 13).Value = -35
$ws.Cells.Item(4, 14).Value = -347.333336
$ws.Cells.Item(75, 8).Value = 12362.5
$ws.Cells.Item(75, 9).Value = 4780
$ws.Cells.Item(75, 10).Value = 25000
$ws.Cells.Item(75, 11).Value = 4780
$ws.Cells.Item(75, 12).Value = 25000
$ws.Cells.Item(75, 13).Value = -3844
$ws.Cells.Item(75, 14).Value = -26872
$ws.Cells.Item(78, 8).Value = 12362.5
$ws.Cells.Item(78, 9).Value = 4780
$ws.Cells.Item(78, 10).Value = 25000
$ws.Cells.Item(78, 11).Value = 14340
$ws.Cells.Item(78, 12).Value = 75000
$ws.Cells.Item(78, 13).Value = -9660
$ws.Cells.Item(78, 14).Value = -84360
$ws.Cells.Item(86, 8).Value = 44741.582
$ws.Cells.Item(86, 10).Value = 103219.8
$ws.Cells.Item(86, 12).Value = 103219.8
$ws.Cells.Item(86, 14).Value = -105465.8
$ws.Cells.Item(89, 8).Value = 44741.582
$ws.Cells.Item(89, 10).Value = 103219.8
$ws.Cells.Item(89, 12).Value = 516099
$ws.Cells.Item(89, 14).Value = -527331
$ws.Cells.Item(94, 8).Value = 999.9091
$ws.Cells.Item(94, 9).Value = 749.875
$ws.Cells.Item(94, 10).Value = 1666.6666
$ws.Cells.Item(94, 11).Value = 749.875
$ws.Cells.Item(94, 12).Value = 1666.6666
$ws.Cells.Item(94, 13).Value = -298.875
$ws.Cells.Item(94, 14).Value = -2568.6666
$ws.Cells.Item(134, 8).Value = 5412.2666
$ws.Cells.Item(134, 9).Value = 5172.3706
$ws.Cells.Item(134, 11).Value = 15517.1118
$ws.Cells.Item(134, 13).Value = -12982.1118

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(10, 8).Value = 22836.564
$ws.Cells.Item(10, 10).Value = 64505
$ws.Cells.Item(10, 12).Value = 64505
$ws.Cells.Item(10, 14).Value = -64783
$ws.Cells.Item(41, 8).Value = 7365.8335
$ws.Cells.Item(41, 9).Value = 2732.6667
$ws.Cells.Item(41, 10).Value = 11999
$ws.Cells.Item(41, 11).Value = 2732.6667
$ws.Cells.Item(41, 12).Value = 11999
$ws.Cells.Item(41, 13).Value = -2304.6667
$ws.Cells.Item(41, 14).Value = -12855
$ws.Cells.Item(50, 8).Value = 9624
$ws.Cells.Item(50, 10).Value = 9624
$ws.Cells.Item(50, 12).Value = 9624
$ws.Cells.Item(50, 14).Value = -10874
$ws.Cells.Item(51, 8).Value = 23082
$ws.Cells.Item(51, 10).Value = 23082
$ws.Cells.Item(51, 12).Value = 23082
$ws.Cells.Item(51, 14).Value = -24554
$ws.Cells.Item(59, 8).Value = 0
$ws.Cells.Item(59, 10).Value = 0
$ws.Cells.Item(59, 12).Value = 0
$ws.Cells.Item(59, 14).ClearContents()
$ws.Cells.Item(60, 8).Value = 12300
$ws.Cells.Item(60, 10).Value = 0
$ws.Cells.Item(60, 12).Value = 0
$ws.Cells.Item(60, 14).ClearContents()
$ws.Cells.Item(61, 8).Value = 23082
$ws.Cells.Item(61, 10).Value = 23082
$ws.Cells.Item(61, 12).Value = 23082
$ws.Cells.Item(61, 14).Value = -23778
$ws.Cells.Item(68, 8).Value = 29947.5
$ws.Cells.Item(68, 10).Value = 29947.5
$ws.Cells.Item(68, 12).Value = 29947.5
$ws.Cells.Item(68, 14).Value = -31445.5
$ws.Cells.Item(71, 8).Value = 29947.5
$ws.Cells.Item(71, 10).Value = 29947.5
$ws.Cells.Item(71, 12).Value = 89842.5
$ws.Cells.Item(71, 14).Value = -97330.5
$ws.Cells.Item(74, 8).Value = 23212.666
$ws.Cells.Item(74, 10).Value = 23212.666
$ws.Cells.Item(74, 12).Value = 23212.666
$ws.Cells.Item(74, 14).Value = -24960.666
$ws.Cells.Item(77, 8).Value = 23212.666
$ws.Cells.Item(77, 10).Value = 23212.666
$ws.Cells.Item(77, 12).Value = 69637.99800000001
$ws.Cells.Item(77, 14).Value = -78373.99800000001
$ws.Cells.Item(105, 8).Value = 4078.182
$ws.Cells.Item(105, 10).Value = 2981.25
$ws.Cells.Item(105, 12).Value = 2981.25
$ws.Cells.Item(105, 14).Value = -6475.25
$ws.Cells.Item(122, 8).Value = 2770.4167
$ws.Cells.Item(122, 9).Value = 2288.7368
$ws.Cells.Item(122, 10).Value = 4600.8
$ws.Cells.Item(122, 11).Value = 6866.2104
$ws.Cells.Item(122, 12).Value = 13802.4
$ws.Cells.Item(122, 13).Value = -4416.2104
$ws.Cells.Item(122, 14).Value = -18702.4

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(132, 8).Value = 1657.4445
$ws.Cells.Item(132, 9).Value = 1346.2727
$ws.Cells.Item(132, 10).Value = 1794.36
$ws.Cells.Item(132, 11).Value = 12116.4543
$ws.Cells.Item(132, 12).Value = 16149.24
$ws.Cells.Item(132, 13).Value = -9586.454299999999
$ws.Cells.Item(132, 14).Value = -21209.24
$ws.Cells.Item(140, 8).Value = 6975.6
$ws.Cells.Item(140, 9).Value = 1400
$ws.Cells.Item(140, 11).Value = 4200
$ws.Cells.Item(140, 13).Value = 980
$ws.Cells.Item(141, 8).Value = 4132.5
$ws.Cells.Item(141, 9).Value = 765
$ws.Cells.Item(141, 10).Value = 7500
$ws.Cells.Item(141, 11).Value = 2295
$ws.Cells.Item(141, 12).Value = 22500
$ws.Cells.Item(141, 13).Value = 2885
$ws.Cells.Item(141, 14).Value = -32860

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(82, 8).Value = 21500
$ws.Cells.Item(82, 9).Value = 13000
$ws.Cells.Item(82, 11).Value = 13000
$ws.Cells.Item(82, 13).Value = -12617
$ws.Cells.Item(85, 8).Value = 21500
$ws.Cells.Item(85, 9).Value = 13000
$ws.Cells.Item(85, 11).Value = 13000
$ws.Cells.Item(85, 13).Value = -11674
$ws.Cells.Item(122, 8).Value = 4907.048
$ws.Cells.Item(122, 9).Value = 3818.4443
$ws.Cells.Item(122, 10).Value = 6866.533
$ws.Cells.Item(122, 11).Value = 11455.3329
$ws.Cells.Item(122, 12).Value = 20599.599
$ws.Cells.Item(122, 13).Value = -9005.332900000001
$ws.Cells.Item(122, 14).Value = -25499.599
$ws.Cells.Item(126, 8).Value = 717915.2
$ws.Cells.Item(126, 9).Value = 2299.8333
$ws.Cells.Item(126, 10).Value = 1254626.8
$ws.Cells.Item(126, 11).Value = 6899.499899999999
$ws.Cells.Item(126, 12).Value = 3763880.4
$ws.Cells.Item(126, 13).Value = -4429.499899999999
$ws.Cells.Item(126, 14).Value = -3768820.4
$ws.Cells.Item(132, 8).Value = 2862.2034
$ws.Cells.Item(132, 9).Value = 2656.1538
$ws.Cells.Item(132, 10).Value = 3264
$ws.Cells.Item(132, 11).Value = 7968.4614
$ws.Cells.Item(132, 12).Value = 9792
$ws.Cells.Item(132, 13).Value = -5438.4614
$ws.Cells.Item(132, 14).Value = -14852

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(14, 8).Value = 167944.7
$ws.Cells.Item(14, 9).Value = 1251050
$ws.Cells.Item(14, 10).Value = 23530.666
$ws.Cells.Item(14, 11).Value = 1251050
$ws.Cells.Item(14, 12).Value = 23530.666
$ws.Cells.Item(14, 13).Value = -1250878
$ws.Cells.Item(14, 14).Value = -23874.666
$ws.Cells.Item(132, 8).Value = 2174.6
$ws.Cells.Item(132, 9).Value = 1495.6765
$ws.Cells.Item(132, 10).Value = 3273.8096
$ws.Cells.Item(132, 11).Value = 4487.029500000001
$ws.Cells.Item(132, 12).Value = 9821.4288
$ws.Cells.Item(132, 13).Value = -1957.029500000001
$ws.Cells.Item(132, 14).Value = -14881.4288
$ws.Cells.Item(140, 8).Value = 51342.855
$ws.Cells.Item(140, 10).Value = 51342.855
$ws.Cells.Item(140, 12).Value = 51342.855
$ws.Cells.Item(140, 14).Value = -61702.855

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 12).Value = 0
$ws.Cells.Item(7, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 12785.813
$ws.Cells.Item(132, 9).Value = 3336.7715
$ws.Cells.Item(132, 10).Value = 54125.375
$ws.Cells.Item(132, 11).Value = 10010.3145
$ws.Cells.Item(132, 12).Value = 162376.125
$ws.Cells.Item(132, 13).Value = -7480.3145
$ws.Cells.Item(132, 14).Value = -167436.125
$ws.Cells.Item(136, 8).Value = 3162.7036
$ws.Cells.Item(136, 9).Value = 3005.8
$ws.Cells.Item(136, 10).Value = 3358.8333
$ws.Cells.Item(136, 11).Value = 9017.400000000001
$ws.Cells.Item(136, 12).Value = 10076.4999
$ws.Cells.Item(136, 13).Value = -3005.8
$ws.Cells.Item(136, 14).Value = -15176.4999
